# Dev Guide change for the Logic component sequence diagram.
#
# 1) The "datetimeFigureOut" date placeholders on the slide master, every
#    custom (slide) layout, and the notes master were re-cached from
#    "2/6/2017" to "3/23/17" (the deck was re-saved on a later date).
# 2) The "deletePerson(p)" call-out textbox on the slide was renamed to
#    "deleteTask(p)", with the trailing "(p)" run split into "(p" and ")"
#    (same formatting) so it matches the author's re-typed run boundary.

$p = $ppt.ActivePresentation

$oldDate = "2/6/2017"
$newDate = "3/23/17"

function Update-DateShape($sh) {
    if ($sh.HasTextFrame) {
        $tr = $sh.TextFrame.TextRange
        if ($tr.Text -eq $oldDate) {
            $tr.Text = $newDate
        }
    }
}

# --- Slide master ---
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    Update-DateShape($master.Shapes.Item($i))
}

# --- Every slide (custom) layout ---
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        Update-DateShape($layout.Shapes.Item($i))
    }
}

# --- Notes master ---
$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    Update-DateShape($notesMaster.Shapes.Item($i))
}

# --- "deletePerson(p)" -> "deleteTask(p)" textbox on slide 1 ---
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $sh = $slide.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        $tr = $sh.TextFrame.TextRange
        if ($tr.Text -eq "deletePerson(p)") {
            # Rename the method, keeping "(p)" as-is for now (its own run).
            $nameRange = $tr.Characters(1, 12)
            $nameRange.Text = "deleteTask"

            # Re-text changed "deletePerson" (13 chars) -> "deleteTask" (10
            # chars); the full string is now "deleteTask(p)" (14 chars).
            # Split the trailing "(p)" run into "(p" + ")" by nudging the
            # font color on just the "(p" slice (identical color), which
            # forces a run break without altering the visible formatting.
            $parenRange = $tr.Characters(11, 2)
            $parenRange.Font.Color.RGB = $parenRange.Font.Color.RGB
        }
    }
}
